$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 0.7
$ws.Range("B17").Value = 4882798.597252
$ws.Range("C17").Value = 45119.137866
$ws.Range("D17").Value = 4837679.459386
$ws.Range("E17").Value = 38955.36240733333
$ws.Range("F17").Value = 1447410.006967
$ws.Range("G17").Value = 20641.353272
$ws.Range("H17").Value = 1426768.653694
$ws.Range("I17").Value = 41516.32374566666
